# lab_07-state_management.pptx — "fixed typos in lesson 7"
#
# 1) Turn on the "slide number" footer placeholder for every content slide
#    (slides 2..N — the title slide, slide 1, is left untouched) via the
#    standard Header & Footer dialog equivalent (HeadersFooters.SlideNumber).
#    PowerPoint names the freshly-minted placeholder shape
#    "Slide Number Placeholder <Id-1>", so we rename it to match.
# 2) Fix a couple of wording typos on the "State" slide and the
#    "ChangeNotifier" slide.

$p = $ppt.ActivePresentation

for ($i = 2; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $s.HeadersFooters.SlideNumber.Visible = 1
    $shp = $s.Shapes.Item($s.Shapes.Count)
    $shp.Name = "Slide Number Placeholder " + ($shp.Id - 1)
}

# --- "State" slide: wording tweaks -----------------------------------
$slideState = $p.Slides.Item(5)
$bodyState = $slideState.Shapes.Item(2).TextFrame.TextRange

$introRun = $bodyState.Paragraphs(1).Runs(1)
$introRun.Text = "State stands for everything that is necessary to define how the app and its screen behave and look at some point in time:"

$ephemeralRun = $bodyState.Paragraphs(8).Runs(2)
$localStart = $ephemeralRun.Start + $ephemeralRun.Text.IndexOf("local state")
$bodyState.Characters($localStart, "local state".Length).Font.Italic = 1

$appStateRun = $bodyState.Paragraphs(9).Runs(2)
$sharedStart = $appStateRun.Start + $appStateRun.Text.IndexOf("shared state")
$bodyState.Characters($sharedStart, "shared state".Length).Font.Italic = 1

# --- "ChangeNotifier" slide: wording tweaks --------------------------
$slideNotifier = $p.Slides.Item(16)
$bodyNotifier = $slideNotifier.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1)

$bodyNotifier.Runs(2).Text = " is a class that can notify "
$bodyNotifier.Runs(3).Text = "listeners "
$bodyNotifier.Runs(4).Text = "of any changes in the state. "
